$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$ws2 = $wb.Worksheets.Item("Lương")

# ================= Sheet 1: "Đơn sale chính" =================
# Row 1: headers (all text)
$ws1.Range("A1").Value = "Tiền tố"
$ws1.Range("B1").Value = "Mã dịch vụ"
$ws1.Range("C1").Value = "Ngày thực hiện"
$ws1.Range("D1").Value = "Cơ sở"
$ws1.Range("E1").Value = "Khách hàng"
$ws1.Range("F1").Value = "Nguồn khách"
$ws1.Range("G1").Value = "Nhóm dịch vụ"
$ws1.Range("H1").Value = "Tên dịch vụ"
$ws1.Range("I1").Value = "Sale chính"
$ws1.Range("J1").Value = "Đơn giá gốc"
$ws1.Range("K1").Value = "Sale phụ"
$ws1.Range("L1").Value = "Upsale"
$ws1.Range("M1").Value = "Đơn giá"
$ws1.Range("N1").Value = "Thanh toán lần đầu"
$ws1.Range("O1").Value = "Trả sau"
$ws1.Range("P1").Value = "Đã thanh toán"
$ws1.Range("Q1").Value = "Dư nợ"
$ws1.Range("R1").Value = "Bác sĩ 1"
$ws1.Range("S1").Value = "Bác sĩ 2"
$ws1.Range("T1").Value = "Phụ phẫu 1"
$ws1.Range("U1").Value = "Phụ phẫu 2"
$ws1.Range("V1").Value = "Công phụ phẫu 1"
$ws1.Range("W1").Value = "Công phụ phẫu 2"
$ws1.Range("X1").Value = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("Y1").Value = "Tỉ lệ chiết khấu sale phụ"
$ws1.Range("Z1").Value = "Chiết khấu sale chính"
$ws1.Range("AA1").Value = "Chiết khấu sale phụ"

# Row 2: order detail data
$ws1.Range("A2").Value = "HD-LUXURY"
$ws1.Range("B2").Value = 555
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "07-16-2024"
$ws1.Range("D2").Value = "LONG XUYÊN"
$ws1.Range("E2").Value = "Nguyễn thị hồng cúc"
$ws1.Range("F2").Value = "Cá nhân"
$ws1.Range("G2").Value = "Tiêm"
$ws1.Range("H2").Value = "Tiêm Filler"
$ws1.Range("I2").Value = "Lê Văn Linh"
$ws1.Range("J2").Value = 7000000
$ws1.Range("M2").Value = 7000000
$ws1.Range("N2").Value = 7000000
$ws1.Range("O2").Value = 0
$ws1.Range("P2").Value = 7000000
$ws1.Range("Q2").Value = 0
$ws1.Range("R2").Value = "Đặng Ngọc Mai"
$ws1.Range("T2").Value = "Sang sang"
$ws1.Range("V2").Value = 50000
$ws1.Range("W2").Value = 0
$ws1.Range("X2").Value = 0.13
$ws1.Range("Y2").Value = 0
$ws1.Range("Z2").Value = 910000
$ws1.Range("AA2").Value = 0

# Row 3: totals row
$ws1.Range("A3").Value = "Tổng"
$ws1.Range("B3").Value = 1
$ws1.Range("J3").Value = 7000000
$ws1.Range("L3").Value = 0
$ws1.Range("M3").Value = 7000000
$ws1.Range("N3").Value = 7000000
$ws1.Range("O3").Value = 0
$ws1.Range("P3").Value = 7000000
$ws1.Range("Q3").Value = 0
$ws1.Range("V3").Value = 50000
$ws1.Range("W3").Value = 0
$ws1.Range("X3").Value = 0.13
$ws1.Range("Y3").Value = 0
$ws1.Range("Z3").Value = 910000
$ws1.Range("AA3").Value = 0


# ================= Sheet 2: "Lương" =================
# Row 1-3: updated values; Row 4: unchanged; Rows 5-10: unchanged (already correct)
$ws2.Range("B2").Value = 16
$ws2.Range("B3").Value = 560000

# Rows 11-31: labels shift down by one (new "Ứng lương" rows inserted) + new totals block
$ws2.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$ws2.Range("B11").Value = -2249000
$ws2.Range("A12").Value = "Lương cơ bản tại LONG XUYÊN"
$ws2.Range("B12").Value = 6571428.571428572
$ws2.Range("A13").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws2.Range("B13").Value = 910000
$ws2.Range("A14").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws2.Range("B14").Value = 0
$ws2.Range("A15").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws2.Range("B15").Value = 0
$ws2.Range("A16").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws2.Range("B16").Value = 0
$ws2.Range("A17").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws2.Range("B17").Value = 0
$ws2.Range("A18").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws2.Range("B18").Value = 0
$ws2.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$ws2.Range("B19").Value = -0.0
$ws2.Range("A20").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws2.Range("B20").Value = 9857142.857142856
$ws2.Range("A21").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws2.Range("B21").Value = 0
$ws2.Range("A22").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws2.Range("B22").Value = 0
$ws2.Range("A23").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws2.Range("B23").Value = 0
$ws2.Range("A24").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws2.Range("B24").Value = 0
$ws2.Range("A25").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws2.Range("B25").Value = 0
$ws2.Range("A26").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws2.Range("B26").Value = 0
$ws2.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$ws2.Range("B27").Value = -0.0
$ws2.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$ws2.Range("B28").Value = 2511000
$ws2.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$ws2.Range("B29").Value = 7481428.571428572
$ws2.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$ws2.Range("B30").Value = 9857142.857142856
$ws2.Range("A31").Value = "Tổng lương"
$ws2.Range("B31").Value = 19849571.42857143

